# mock-forecast: rename disbursement categories to Purchases / Taxes / Transfers,
# and fold the old "Disbursement 4" row into the (renamed) "Transfers" row.

$wb = $excel.ActiveWorkbook

# --- "data" sheet: rename categories B3:B5, remove the 4th category row ---
$data = $wb.Worksheets.Item("data")

$data.Range("B3").Value = "Purchases"
$data.Range("B4").Value = "Taxes"
$data.Range("B5").Value = "Transfers"

# Row 6 ("Disbursement 4") is no longer a separate category; its numbers were
# merged into row 5 on the "mpc" sheet, so clear it here.
$data.Range("B6:K6").ClearContents()

$data.Activate()
$data.Range("B3").Select()

# --- "mpc" sheet: move the old row-6 percentages up onto row 5 (now "Transfers") ---
$mpc = $wb.Worksheets.Item("mpc")

$mpc.Range("C5").Value = 0.3
$mpc.Range("D5").Value = 0.2
$mpc.Range("E5").Value = 0.1
$mpc.Range("F5").Value = 0.05
$mpc.Range("G5").Value = 0.05
$mpc.Range("H5").Value = 0.02
$mpc.Range("I5").Value = 0.02
$mpc.Range("J5").Value = 0.02
$mpc.Range("K5").Value = 0.02
$mpc.Range("L5").Value = 0.02

# Old row 6 numbers have moved to row 5, so clear them from row 6.
$mpc.Range("C6:L6").ClearContents()

$mpc.Activate()
$mpc.Range("I23").Select()
